$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 27778714
$ws.Range("I33").Value = 41667690
$ws.Range("K33").Value = 41667690
$ws.Range("M33").Value = -41667461
$ws.Range("H86").Value = 2082.7222
$ws.Range("I86").Value = 2295
$ws.Range("J86").Value = 1817.375
$ws.Range("K86").Value = 2295
$ws.Range("L86").Value = 1817.375
$ws.Range("M86").Value = -1172
$ws.Range("N86").Value = -4063.375
$ws.Range("H89").Value = 2082.7222
$ws.Range("I89").Value = 2295
$ws.Range("J89").Value = 1817.375
$ws.Range("K89").Value = 11475
$ws.Range("L89").Value = 9086.875
$ws.Range("M89").Value = -5859
$ws.Range("N89").Value = -20318.875
$ws.Range("H98").Value = 1864.7142
$ws.Range("I98").Value = 1864.7142
$ws.Range("K98").Value = 1864.7142
$ws.Range("M98").Value = -366.7141999999999
$ws.Range("H116").Value = 4156.1113
$ws.Range("I116").Value = 3201
$ws.Range("K116").Value = 3201
$ws.Range("M116").Value = 241
$ws.Range("H122").Value = 1864.7142
$ws.Range("I122").Value = 1864.7142
$ws.Range("K122").Value = 5594.142599999999
$ws.Range("M122").Value = -3144.142599999999
$ws.Range("H137").Value = 3953.7097
$ws.Range("I137").Value = 4056.25
$ws.Range("K137").Value = 12168.75
$ws.Range("M137").Value = -9618.75
$ws.Range("H138").Value = 2362.2856
$ws.Range("I138").Value = 1159.6364
$ws.Range("J138").Value = 3342.2222
$ws.Range("K138").Value = 3478.9092
$ws.Range("L138").Value = 10026.6666
$ws.Range("M138").Value = 1661.0908
$ws.Range("N138").Value = -20306.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7714.8438
$ws.Range("I32").Value = 5976.6294
$ws.Range("K32").Value = 5976.6294
$ws.Range("M32").Value = -5689.6294
$ws.Range("H45").Value = 206099.9
$ws.Range("I45").Value = 290142.72
$ws.Range("K45").Value = 290142.72
$ws.Range("M45").Value = -289765.72
$ws.Range("H61").Value = 5884.452
$ws.Range("I61").Value = 5998.85
$ws.Range("K61").Value = 5998.85
$ws.Range("M61").Value = -5786.85
$ws.Range("H74").Value = 4542.55
$ws.Range("I74").Value = 3380.611
$ws.Range("K74").Value = 3380.611
$ws.Range("M74").Value = -2506.611
$ws.Range("H77").Value = 4542.55
$ws.Range("I77").Value = 3380.611
$ws.Range("K77").Value = 16903.055
$ws.Range("M77").Value = -12535.055
$ws.Range("H132").Value = 3161.1
$ws.Range("I132").Value = 2901.375
$ws.Range("J132").Value = 4200
$ws.Range("K132").Value = 8704.125
$ws.Range("L132").Value = 12600
$ws.Range("M132").Value = -6174.125
$ws.Range("N132").Value = -17660
$ws.Range("H136").Value = 5884.452
$ws.Range("I136").Value = 5998.85
$ws.Range("K136").Value = 17996.55
$ws.Range("M136").Value = -15446.55

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2912.087
$ws.Range("I20").Value = 2593.4
$ws.Range("J20").Value = 3509.625
$ws.Range("K20").Value = 2593.4
$ws.Range("L20").Value = 3509.625
$ws.Range("M20").Value = -2346.4
$ws.Range("N20").Value = -4003.625
$ws.Range("H97").Value = 14389.667
$ws.Range("I97").Value = 9849
$ws.Range("K97").Value = 9849
$ws.Range("M97").Value = -8858
$ws.Range("H134").Value = 6161.609
$ws.Range("I134").Value = 5998.077
$ws.Range("K134").Value = 17994.231
$ws.Range("M134").Value = -15459.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5486.7144
$ws.Range("I31").Value = 3283.1667
$ws.Range("K31").Value = 3283.1667
$ws.Range("M31").Value = -2988.1667
$ws.Range("H34").Value = 5486.7144
$ws.Range("I34").Value = 3283.1667
$ws.Range("K34").Value = 3283.1667
$ws.Range("M34").Value = -3081.1667
$ws.Range("H62").Value = 4658.2
$ws.Range("I62").Value = 4250
$ws.Range("J62").Value = 4930.3335
$ws.Range("K62").Value = 4250
$ws.Range("L62").Value = 4930.3335
$ws.Range("M62").Value = -3626
$ws.Range("N62").Value = -6178.3335
$ws.Range("H65").Value = 4658.2
$ws.Range("I65").Value = 4250
$ws.Range("J65").Value = 4930.3335
$ws.Range("K65").Value = 21250
$ws.Range("L65").Value = 24651.6675
$ws.Range("M65").Value = -18130
$ws.Range("N65").Value = -30891.6675
$ws.Range("H97").Value = 28433.8
$ws.Range("J97").Value = 28433.8
$ws.Range("L97").Value = 28433.8
$ws.Range("N97").Value = -30415.8
$ws.Range("H99").Value = 8124.25
$ws.Range("I99").Value = 8124.25
$ws.Range("K99").Value = 8124.25
$ws.Range("M99").Value = -6626.25
$ws.Range("H105").Value = 71429930
$ws.Range("J105").Value = 1366.3334
$ws.Range("L105").Value = 1366.3334
$ws.Range("N105").Value = -4860.3334
$ws.Range("H126").Value = 8124.25
$ws.Range("I126").Value = 8124.25
$ws.Range("K126").Value = 24372.75
$ws.Range("M126").Value = -21902.75
$ws.Range("H132").Value = 7154.3335
$ws.Range("I132").Value = 5157.8
$ws.Range("K132").Value = 15473.4
$ws.Range("M132").Value = -12943.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 241.5
$ws.Range("I44").Value = 239.8
$ws.Range("J44").Value = 242.71428
$ws.Range("K44").Value = 719.4000000000001
$ws.Range("L44").Value = 728.14284
$ws.Range("M44").Value = -321.4000000000001
$ws.Range("N44").Value = -1524.14284
$ws.Range("H134").Value = 1175.0555
$ws.Range("I134").Value = 1175.0555
$ws.Range("K134").Value = 3525.1665
$ws.Range("M134").Value = 1544.8335
$ws.Range("H137").Value = 5009
$ws.Range("I137").Value = 3867.7
$ws.Range("K137").Value = 11603.1
$ws.Range("M137").Value = -6503.099999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 62001.5
$ws.Range("J121").Value = 62001.5
$ws.Range("L121").Value = 62001.5
$ws.Range("N121").Value = -65495.5
$ws.Range("H132").Value = 6279.75
$ws.Range("I132").Value = 5454.8335
$ws.Range("K132").Value = 16364.5005
$ws.Range("M132").Value = -13834.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4872.2
$ws.Range("I40").Value = 4846.5
$ws.Range("K40").Value = 4846.5
$ws.Range("M40").Value = -4710.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 7070707
$ws.Range("J5").Value = 7070707
$ws.Range("L5").Value = 7070707
$ws.Range("N5").Value = -7070931
$ws.Range("H136").Value = 8851.117
$ws.Range("I136").Value = 8031.3335
$ws.Range("K136").Value = 24094.0005
$ws.Range("M136").Value = -21544.0005
$ws.Range("H140").Value = 89794.55
$ws.Range("J140").Value = 89794.55
$ws.Range("L140").Value = 89794.55
$ws.Range("N140").Value = -100154.55
$ws.Range("H141").Value = 80623.5
$ws.Range("J141").Value = 80623.5
$ws.Range("L141").Value = 80623.5
$ws.Range("N141").Value = -90983.5
